$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v0")

# Row 7: v0-protocol-4v0.json / LC - short waits for demo
$ws.Range("A7").Value = "v0-protocol-4v0.json"
$ws.Range("B7").Value = "LC - short waits for demo"
$ws.Range("C7").Value = "5 mL"
$ws.Range("D7").Value = "5 mL"
$ws.Range("E7").Value = "15 sec"
$ws.Range("F7").Value = "0.5 mL"
$ws.Range("G7").Value = "15 mL/hr"
$ws.Range("H7").Value = "15 mL/hr"
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = "15 sec"
$ws.Range("K7").Value = "N"

# Row 8: v0-protocol-5v0.json / PANC 1.5 mL Plasma (Yang CCR 2020)
$ws.Range("B8").Value = "PANC 1.5 mL Plasma (Yang CCR 2020)"
$ws.Range("C8").Value = "5 mL"
$ws.Range("D8").Value = "5 mL"
$ws.Range("E8").Value = "1 hour"
$ws.Range("F8").Value = "1.5 mL"
$ws.Range("G8").Value = "1.5 mL/hr"
$ws.Range("H8").Value = "1.5 mL/hr"
$ws.Range("H8").Font.Color = 255
$ws.Range("I8").Value = 1000
$ws.Range("I8").HorizontalAlignment = -4131
$ws.Range("J8").Value = "3 mins (QIAZOL - 700 uL)"
$ws.Range("J8").Interior.Color = 65535
$ws.Range("K8").Value = "N"
$ws.Range("L8").Value = "Pull 700 uL to WASTE SYRINGE, 3 min incubation, then final pull to lysate syringe"
$ws.Range("L8").Interior.Color = 65535
$ws.Range("A8").Value = "v0-protocol-5v0.json"

# Update selection/view to reflect final edit state
$ws.Range("A8").Select()
